$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values between row 2 and row 4 for columns D, J, K, L, M, O, P
$row2_D = $ws.Range("D2").Value2
$row2_J = $ws.Range("J2").Value2
$row2_K = $ws.Range("K2").Value2
$row2_L = $ws.Range("L2").Value2
$row2_M = $ws.Range("M2").Value2
$row2_O = $ws.Range("O2").Value2
$row2_P = $ws.Range("P2").Value2

$row4_D = $ws.Range("D4").Value2
$row4_J = $ws.Range("J4").Value2
$row4_K = $ws.Range("K4").Value2
$row4_L = $ws.Range("L4").Value2
$row4_M = $ws.Range("M4").Value2
$row4_O = $ws.Range("O4").Value2
$row4_P = $ws.Range("P4").Value2

$ws.Range("D2").Value2 = $row4_D
$ws.Range("J2").Value2 = $row4_J
$ws.Range("K2").Value2 = $row4_K
$ws.Range("L2").Value2 = $row4_L
$ws.Range("M2").Value2 = $row4_M
$ws.Range("O2").Value2 = $row4_O
$ws.Range("P2").Value2 = $row4_P

$ws.Range("D4").Value2 = $row2_D
$ws.Range("J4").Value2 = $row2_J
$ws.Range("K4").Value2 = $row2_K
$ws.Range("L4").Value2 = $row2_L
$ws.Range("M4").Value2 = $row2_M
$ws.Range("O4").Value2 = $row2_O
$ws.Range("P4").Value2 = $row2_P
